# Insert a new weekly record for "Haba" at Vega Monumental Concepción.
# This shifts the existing rows 5-15 down to 6-16 and fills row 5 with
# the new data point (date 2021-10-29, volumen 220, etc.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 5..15 down by one row (insert a blank row before row 5).
$ws.Rows.Item(5).Insert()

# Populate the new row 5 with the new record's values.
$ws.Cells.Item(5, 1).Value = 11
$ws.Cells.Item(5, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(5, 3).Value = "Bíobío"
$ws.Cells.Item(5, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(5, 4).Value = 44498
$ws.Cells.Item(5, 5).Value = 8
$ws.Cells.Item(5, 6).Value = 100112026
$ws.Cells.Item(5, 7).Value = "Haba"
$ws.Cells.Item(5, 8).Value = "Sin especificar"
$ws.Cells.Item(5, 9).Value = "Primera"
$ws.Cells.Item(5, 10).Value = 220
$ws.Cells.Item(5, 11).Value = 7000
$ws.Cells.Item(5, 12).Value = 7500
$ws.Cells.Item(5, 13).Value = 7273
$ws.Cells.Item(5, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(5, 15).Value = "Región Metropolitana"
$ws.Cells.Item(5, 16).Value = 291
$ws.Cells.Item(5, 17).Value = 25
$ws.Cells.Item(5, 18).Value = "Hortaliza"
